$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D retains its Text format so numeric-looking values
# (e.g. "1.001", "313.57") stay as text instead of being converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "24.646.36"
$ws.Cells.Item(2, 5).Value = "  -1.24%  "
$ws.Cells.Item(3, 4).Value = "1.673.91"
$ws.Cells.Item(3, 5).Value = "  -2.18%  "
$ws.Cells.Item(4, 4).Value = "1.001"
$ws.Cells.Item(4, 5).Value = "  -0.12%  "
$ws.Cells.Item(5, 4).Value = "313.57"
$ws.Cells.Item(5, 5).Value = "  -0.71%  "
$ws.Cells.Item(6, 5).Value = "  -0.19%  "
$ws.Cells.Item(7, 4).Value = "0.3938"
$ws.Cells.Item(7, 5).Value = "  -2.42%  "
$ws.Cells.Item(8, 4).Value = "0.3946"
$ws.Cells.Item(8, 5).Value = "  -3.07%  "
$ws.Cells.Item(9, 4).Value = "1.001"
$ws.Cells.Item(9, 5).Value = "  -0.09%  "
$ws.Cells.Item(10, 4).Value = "1.400"
$ws.Cells.Item(10, 5).Value = "  -5.53%  "
$ws.Cells.Item(11, 4).Value = "50.96"
$ws.Cells.Item(11, 5).Value = "  -5.36%  "
$ws.Cells.Item(12, 4).Value = "0.08644"
$ws.Cells.Item(12, 5).Value = "  -2.10%  "
$ws.Cells.Item(13, 4).Value = "25.27"
$ws.Cells.Item(13, 5).Value = "  -4.55%  "
$ws.Cells.Item(14, 4).Value = "7.323"
$ws.Cells.Item(14, 5).Value = "  -2.68%  "
$ws.Cells.Item(15, 4).Value = "0.00001317"
$ws.Cells.Item(15, 5).Value = "  -2.13%  "
$ws.Cells.Item(16, 4).Value = "7.697"
$ws.Cells.Item(16, 5).Value = "  -4.27%  "
$ws.Cells.Item(17, 4).Value = "1.680.03"
$ws.Cells.Item(17, 5).Value = "  -1.10%  "
$ws.Cells.Item(18, 4).Value = "93.96"
$ws.Cells.Item(18, 5).Value = "  -1.33%  "
$ws.Cells.Item(19, 4).Value = "0.07016"
$ws.Cells.Item(19, 5).Value = "  -2.40%  "
$ws.Cells.Item(20, 5).Value = "  +0.47%  "
$ws.Cells.Item(21, 4).Value = "7.086"
$ws.Cells.Item(21, 5).Value = "  -2.71%  "
$ws.Cells.Item(22, 5).Value = "  -0.36%  "
$ws.Cells.Item(23, 4).Value = "13.92"
$ws.Cells.Item(23, 5).Value = "  -4.26%  "
$ws.Cells.Item(24, 4).Value = "24.649.98"
$ws.Cells.Item(24, 5).Value = "  -1.21%  "
$ws.Cells.Item(25, 4).Value = "2.346"
$ws.Cells.Item(25, 5).Value = "  +0.42%  "
$ws.Cells.Item(26, 4).Value = "2.780"
$ws.Cells.Item(26, 5).Value = "  -4.15%  "
$ws.Cells.Item(27, 4).Value = "23.03"
$ws.Cells.Item(27, 5).Value = "  -0.53%  "
$ws.Cells.Item(28, 4).Value = "5.832"
$ws.Cells.Item(28, 5).Value = "  -9.17%  "
$ws.Cells.Item(29, 4).Value = "160.13"
$ws.Cells.Item(29, 5).Value = "  -1.86%  "
$ws.Cells.Item(30, 4).Value = "145.83"
$ws.Cells.Item(30, 5).Value = "  +1.42%  "
$ws.Cells.Item(31, 4).Value = "8.312"
$ws.Cells.Item(31, 5).Value = "  +1.08%  "
$ws.Cells.Item(32, 4).Value = "2.483"
$ws.Cells.Item(32, 5).Value = "  +9.56%  "
$ws.Cells.Item(33, 4).Value = "1.860.03"
$ws.Cells.Item(33, 5).Value = "  -1.32%  "
$ws.Cells.Item(34, 5).Value = "  -3.19%  "
$ws.Cells.Item(35, 4).Value = "0.08264"
$ws.Cells.Item(35, 5).Value = "  -5.83%  "
$ws.Cells.Item(36, 4).Value = "6.963"
$ws.Cells.Item(36, 5).Value = "  -5.10%  "
$ws.Cells.Item(37, 4).Value = "0.2820"
$ws.Cells.Item(37, 5).Value = "  -2.15%  "
$ws.Cells.Item(38, 4).Value = "0.9926"
$ws.Cells.Item(38, 5).Value = "  -4.09%  "
$ws.Cells.Item(39, 4).Value = "0.09592"
$ws.Cells.Item(39, 5).Value = "  +1.25%  "
$ws.Cells.Item(40, 4).Value = "1.517"
$ws.Cells.Item(40, 5).Value = "  +2.65%  "
$ws.Cells.Item(41, 5).Value = "  -5.52%  "
$ws.Cells.Item(42, 4).Value = "0.7912"
$ws.Cells.Item(42, 5).Value = "  -6.68%  "
$ws.Cells.Item(43, 4).Value = "13.50"
$ws.Cells.Item(43, 5).Value = "  -4.68%  "
$ws.Cells.Item(44, 4).Value = "16.66"
$ws.Cells.Item(44, 5).Value = "  -5.32%  "
$ws.Cells.Item(45, 4).Value = "2.567"
$ws.Cells.Item(45, 5).Value = "  -6.05%  "
$ws.Cells.Item(46, 4).Value = "0.7102"
$ws.Cells.Item(46, 5).Value = "  -4.91%  "
$ws.Cells.Item(47, 4).Value = "4.169"
$ws.Cells.Item(47, 5).Value = "  -1.60%  "
$ws.Cells.Item(48, 4).Value = "0.08656"
$ws.Cells.Item(48, 5).Value = "  +2.90%  "
$ws.Cells.Item(49, 5).Value = "  -0.15%  "
$ws.Cells.Item(50, 4).Value = "1.328"
$ws.Cells.Item(50, 5).Value = "  -4.57%  "
$ws.Cells.Item(51, 4).Value = "137.98"
$ws.Cells.Item(51, 5).Value = "  -2.37%  "
